$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the last data row (old row 19) so the table shrinks from 19 to 18 rows
$ws.Rows.Item(19).Delete()

# Row 2
$ws.Range("A2").Value = "BNK"
$ws.Range("B2").Value = "'2024-02-19"
$ws.Range("C2").Value = "비엔케이제2호스팩"
$ws.Range("D2").Value = "BNK"
$ws.Range("E2").Value = "BNK"
$ws.Range("F2").Value = "'2024-02-22"
$ws.Range("G2").Value = "'2024-03-05"
$ws.Range("H2").Value = 8000
$ws.Range("I2").Value = 4000000
$ws.Range("J2").Value = 2000
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 100

# Row 3
$ws.Range("A3").Value = "DB"
$ws.Range("B3").Value = "'2024-02-26"
$ws.Range("C3").Value = "케이엔알시스템"
$ws.Range("D3").Value = "DB, NH"
$ws.Range("E3").Value = "DB, NH"
$ws.Range("F3").Value = "'2024-02-29"
$ws.Range("G3").Value = "'2024-03-07"
$ws.Range("H3").Value = 14202
$ws.Range("I3").Value = 2104000
$ws.Range("J3").Value = 13500
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 50

# Row 4
$ws.Range("A4").Value = "DB"
$ws.Range("B4").Value = "'2024-01-25"
$ws.Range("C4").Value = "스튜디오삼익"
$ws.Range("D4").Value = "DB"
$ws.Range("E4").Value = "DB"
$ws.Range("F4").Value = "'2024-01-30"
$ws.Range("G4").Value = "'2024-02-06"
$ws.Range("H4").Value = 15300
$ws.Range("I4").Value = 850000
$ws.Range("J4").Value = 18000
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 100

# Row 5
$ws.Range("A5").Value = "NH"
$ws.Range("B5").Value = "'2024-02-26"
$ws.Range("C5").Value = "케이엔알시스템"
$ws.Range("D5").Value = "DB, NH"
$ws.Range("E5").Value = "DB, NH"
$ws.Range("F5").Value = "'2024-02-29"
$ws.Range("G5").Value = "'2024-03-07"
$ws.Range("H5").Value = 14202
$ws.Range("I5").Value = 2104000
$ws.Range("J5").Value = 13500
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 50

# Row 6
$ws.Range("A6").Value = "NH"
$ws.Range("B6").Value = "'2024-02-13"
$ws.Range("C6").Value = "케이웨더"
$ws.Range("D6").Value = "NH"
$ws.Range("E6").Value = "NH"
$ws.Range("F6").Value = "'2024-02-16"
$ws.Range("G6").Value = "'2024-02-22"
$ws.Range("H6").Value = 7000
$ws.Range("I6").Value = 1000000
$ws.Range("J6").Value = 7000
$ws.Range("K6").Value = 0
$ws.Range("L6").Value = 100

# Row 7
$ws.Range("A7").Value = "NH"
$ws.Range("B7").Value = "'2024-03-04"
$ws.Range("C7").Value = "오상헬스케어"
$ws.Range("D7").Value = "NH"
$ws.Range("E7").Value = "NH"
$ws.Range("F7").Value = "'2024-03-07"
$ws.Range("G7").Value = "'2024-03-13"
$ws.Range("H7").Value = 19800
$ws.Range("I7").Value = 990000
$ws.Range("J7").Value = 20000
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 100

# Row 8
$ws.Range("A8").Value = "NH"
$ws.Range("B8").Value = "'2024-03-14"
$ws.Range("C8").Value = "엔젤로보틱스"
$ws.Range("D8").Value = "NH"
$ws.Range("E8").Value = "NH"
$ws.Range("F8").Value = "'2024-03-19"
$ws.Range("G8").Value = "'2024-03-26"
$ws.Range("H8").Value = 32000
$ws.Range("I8").Value = 1600000
$ws.Range("J8").Value = 20000
$ws.Range("K8").Value = 0
$ws.Range("L8").Value = 100

# Row 9
$ws.Range("A9").Value = "SK"
$ws.Range("B9").Value = "'2024-02-20"
$ws.Range("C9").Value = "SK증권제11호스팩"
$ws.Range("D9").Value = "SK"
$ws.Range("E9").Value = "SK"
$ws.Range("F9").Value = "'2024-02-23"
$ws.Range("G9").Value = "'2024-03-04"
$ws.Range("H9").Value = 8000
$ws.Range("I9").Value = 4000000
$ws.Range("J9").Value = 2000
$ws.Range("K9").Value = 0
$ws.Range("L9").Value = 100

# Row 10
$ws.Range("A10").Value = "신영"
$ws.Range("B10").Value = "'2024-01-25"
$ws.Range("C10").Value = "신영스팩10호"
$ws.Range("D10").Value = "신영"
$ws.Range("E10").Value = "신영"
$ws.Range("F10").Value = "'2024-01-30"
$ws.Range("G10").Value = "'2024-02-06"
$ws.Range("H10").Value = 9150
$ws.Range("I10").Value = 4575000
$ws.Range("J10").Value = 2000
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 100

# Row 11
$ws.Range("A11").Value = "신한"
$ws.Range("B11").Value = "'2024-02-14"
$ws.Range("C11").Value = "에이피알"
$ws.Range("D11").Value = "신한"
$ws.Range("E11").Value = "신한, 하나"
$ws.Range("F11").Value = "'2024-02-19"
$ws.Range("G11").Value = "'2024-02-27"
$ws.Range("H11").Value = 75800
$ws.Range("I11").Value = 379000
$ws.Range("J11").Value = 250000
$ws.Range("K11").Value = 0
$ws.Range("L11").Value = 80

# Row 12
$ws.Range("A12").Value = "유안타"
$ws.Range("B12").Value = "'2024-02-20"
$ws.Range("C12").Value = "유안타제15호스팩"
$ws.Range("D12").Value = "유안타"
$ws.Range("E12").Value = "유안타"
$ws.Range("F12").Value = "'2024-02-23"
$ws.Range("G12").Value = "'2024-02-29"
$ws.Range("H12").Value = 13000
$ws.Range("I12").Value = 6500000
$ws.Range("J12").Value = 2000
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 100

# Row 13
$ws.Range("A13").Value = "유진"
$ws.Range("B13").Value = "'2024-02-19"
$ws.Range("C13").Value = "유진스팩10호"
$ws.Range("D13").Value = "유진"
$ws.Range("E13").Value = "유진"
$ws.Range("F13").Value = "'2024-02-22"
$ws.Range("G13").Value = "'2024-02-29"
$ws.Range("H13").Value = 8000
$ws.Range("I13").Value = 4000000
$ws.Range("J13").Value = 2000
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 100

# Row 14
$ws.Range("A14").Value = "키움"
$ws.Range("B14").Value = "'2024-02-13"
$ws.Range("C14").Value = "코셈"
$ws.Range("D14").Value = "키움"
$ws.Range("E14").Value = "키움"
$ws.Range("F14").Value = "'2024-02-16"
$ws.Range("G14").Value = "'2024-02-23"
$ws.Range("H14").Value = 9600
$ws.Range("I14").Value = 600000
$ws.Range("J14").Value = 16000
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 100

# Row 15
$ws.Range("A15").Value = "하나"
$ws.Range("B15").Value = "'2024-02-22"
$ws.Range("C15").Value = "하나31호스팩"
$ws.Range("D15").Value = "하나"
$ws.Range("E15").Value = "하나"
$ws.Range("F15").Value = "'2024-02-27"
$ws.Range("G15").Value = "'2024-03-05"
$ws.Range("H15").Value = 10000
$ws.Range("I15").Value = 5000000
$ws.Range("J15").Value = 2000
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 100

# Row 16
$ws.Range("A16").Value = "하나"
$ws.Range("B16").Value = "'2024-02-14"
$ws.Range("C16").Value = "에이피알"
$ws.Range("D16").Value = "신한"
$ws.Range("E16").Value = "신한, 하나"
$ws.Range("F16").Value = "'2024-02-19"
$ws.Range("G16").Value = "'2024-02-27"
$ws.Range("H16").Value = 18950
$ws.Range("I16").Value = 379000
$ws.Range("J16").Value = 250000
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 20

# Row 17
$ws.Range("A17").Value = "한국"
$ws.Range("B17").Value = "'2024-03-12"
$ws.Range("C17").Value = "삼현"
$ws.Range("D17").Value = "한국"
$ws.Range("E17").Value = "한국"
$ws.Range("F17").Value = "'2024-03-15"
$ws.Range("G17").Value = "'2024-03-21"
$ws.Range("H17").Value = 60000
$ws.Range("I17").Value = 2000000
$ws.Range("J17").Value = 30000
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 100

# Row 18
$ws.Range("A18").Value = "한화"
$ws.Range("B18").Value = "'2024-02-13"
$ws.Range("C18").Value = "이에이트"
$ws.Range("D18").Value = "한화"
$ws.Range("E18").Value = "한화"
$ws.Range("F18").Value = "'2024-02-16"
$ws.Range("G18").Value = "'2024-02-23"
$ws.Range("H18").Value = 22600
$ws.Range("I18").Value = 1130000
$ws.Range("J18").Value = 20000
$ws.Range("K18").Value = 0
$ws.Range("L18").Value = 100

# Reset number formatting on the date columns so they keep the default (unstyled) cell style
$ws.Range("B2:B18").Style = "Normal"
$ws.Range("F2:F18").Style = "Normal"
$ws.Range("G2:G18").Style = "Normal"
